$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: enter data in "natural" (pre-sort) order so that the shared
# string table is built up in the same sequence the author typed it in.
# Existing students (rows 3,5-12) only need their attendance count (C)
# updated -- their names are unchanged.
$ws.Range("C3").Value2  = 2   # Claudiu Druța
$ws.Range("C5").Value2  = 1   # Sorin Fechete
$ws.Range("C6").Value2  = 2   # Delia Negrea
$ws.Range("C7").Value2  = 2   # Răzvan Baroi
$ws.Range("C8").Value2  = 1   # Andra Agud
$ws.Range("C9").Value2  = 2   # Paul Dobroțchi
$ws.Range("C10").Value2 = 2   # Miriam Bacso
$ws.Range("C11").Value2 = 1   # Denisa Cioban
$ws.Range("C12").Value2 = 1   # Amanda Hajdu

# New students, typed into the first empty rows (13-19), in the order they
# were added.
$ws.Range("B13").Value2 = "Silvia Naghi"
$ws.Range("C13").Value2 = 2

$ws.Range("B14").Value2 = "Raul Andrei"
$ws.Range("C14").Value2 = 1

$ws.Range("B15").Value2 = "Levente Nagy"
$ws.Range("C15").Value2 = 2

$ws.Range("B16").Value2 = "Mark Pop"
$ws.Range("C16").Value2 = 2

$ws.Range("B17").Value2 = "Alessandro Vereș-Pop"
$ws.Range("C17").Value2 = 2

$ws.Range("B18").Value2 = "Luca Șeicaru"
$ws.Range("C18").Value2 = 1

$ws.Range("B19").Value2 = "Attila Bunta"
$ws.Range("C19").Value2 = 1

# Existing student renamed (and recounted).
$ws.Range("B4").Value2  = "Daniela Cionca (Mărie)"
$ws.Range("C4").Value2  = 2

# Two more new students.
$ws.Range("B20").Value2 = "Victor Lazăr"
$ws.Range("C20").Value2 = 1

$ws.Range("B21").Value2 = "Codruț Avram"
$ws.Range("C21").Value2 = 1

# --- Phase 2: rearrange rows 3-21 into alphabetical order by first name,
# matching the workbook's sortState (B3:S21). Cell formatting is tied to the
# row, so we move the logical values rather than the rows themselves.
$names = @(
    "Alessandro Vereș-Pop",
    "Amanda Hajdu",
    "Andra Agud",
    "Attila Bunta",
    "Claudiu Druța",
    "Codruț Avram",
    "Daniela Cionca (Mărie)",
    "Delia Negrea",
    "Denisa Cioban",
    "Levente Nagy",
    "Luca Șeicaru",
    "Mark Pop",
    "Miriam Bacso",
    "Paul Dobroțchi",
    "Raul Andrei",
    "Răzvan Baroi",
    "Silvia Naghi",
    "Sorin Fechete",
    "Victor Lazăr"
)
$counts = @{
    "Alessandro Vereș-Pop"   = 2
    "Amanda Hajdu"           = 1
    "Andra Agud"             = 1
    "Attila Bunta"           = 1
    "Claudiu Druța"          = 2
    "Codruț Avram"           = 1
    "Daniela Cionca (Mărie)" = 2
    "Delia Negrea"           = 2
    "Denisa Cioban"          = 1
    "Levente Nagy"           = 2
    "Luca Șeicaru"           = 1
    "Mark Pop"               = 2
    "Miriam Bacso"           = 2
    "Paul Dobroțchi"         = 2
    "Raul Andrei"            = 1
    "Răzvan Baroi"           = 2
    "Silvia Naghi"           = 2
    "Sorin Fechete"          = 1
    "Victor Lazăr"           = 1
}

$row = 3
foreach ($n in $names) {
    $ws.Cells.Item($row, 2).Value2 = $n
    $ws.Cells.Item($row, 3).Value2 = $counts[$n]
    $row = $row + 1
}

# --- Refresh the worksheet's remembered AutoFilter/sort state (B3:S22,
# sorted by column B) to match the already-sorted data above.
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("B3:B22")) | Out-Null
$sort.SetRange($ws.Range("B3:S22"))
$sort.Header = 0
$sort.Apply()

# --- Misc view bits to mirror the saved sheet view/selection.
$ws.Range("B3:C21").Select()
